$d = $word.ActiveDocument

function Set-ParaText($para, $newText) {
    $r = $para.Range
    $r2 = $d.Range($r.Start, $r.End - 1)
    $r2.Text = $newText
}

# ---------------------------------------------------------------------------
# 1) "Review of work undertaken" bullet list: shift text down one slot and
#    fix the first line's wording, then turn the old bookmark-bearing last
#    bullet into plain text and append two brand-new bullets after it.
# ---------------------------------------------------------------------------

Set-ParaText $d.Paragraphs(8) "Connection issues fixed. "
Set-ParaText $d.Paragraphs(9) "Completion of benchmarking process."

$p10 = $d.Paragraphs(10)
$full10 = $p10.Range
$xml10 = @"
<?xml version='1.0' encoding='UTF-8' standalone='yes'?>
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="Heading2Sturley"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr><w:t>Integrate benchmarking process into GUI.</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="Heading2Sturley"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr><w:t>Implement spoof trading – both backend logic and GUI.</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="Heading2Sturley"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr><w:t>Various bug fixes.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$full10.InsertXML($xml10)

# ---------------------------------------------------------------------------
# 2) "This week X hours ..." -> split into three runs, with the number 32
#    in its own run.
# ---------------------------------------------------------------------------

$hoursPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("This week X hours")) {
        $hoursPara = $p
        break
    }
}
$fullHours = $hoursPara.Range
$xmlHours = @"
<?xml version='1.0' encoding='UTF-8' standalone='yes'?>
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="Heading2Sturley"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve">This week </w:t></w:r>
<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr><w:t>32</w:t></w:r>
<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve"> hours have been spent working on the project.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$fullHours.InsertXML($xmlHours)

# ---------------------------------------------------------------------------
# 3) "Plan of work for the next week" section: fill in the empty bulleted
#    paragraph and append two further bullets, moving the _GoBack bookmark
#    onto the new "OR (pending demo with Marco)" bullet.
# ---------------------------------------------------------------------------

$planPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "" -and $p.Style.NameLocal -eq "Heading 2 Sturley") {
        try {
            if ($p.Range.ListFormat.ListString -ne $null) { }
        } catch {}
    }
}

# Locate the empty numbered paragraph right after "Plan of work for the next week"
$planHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim().StartsWith("Plan of work for the next week")) {
        $planHeading = $p
        break
    }
}
$targetPara = $planHeading.Next()

$fullPlan = $targetPara.Range
$xmlPlan = @"
<?xml version='1.0' encoding='UTF-8' standalone='yes'?>
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="Heading2Sturley"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr><w:t>Implement trading restrictions and issuing trade orders to GDAX endpoint.</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="Heading2Sturley"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr><w:t>OR (pending demo with Marco)</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="Heading2Sturley"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr><w:t>Implement neural network trading predictions.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$fullPlan.InsertXML($xmlPlan)

Write-Output "Edit complete"
